$wb = $excel.ActiveWorkbook

# --- Update the "CODIGO CASO DE USO" (H column) values on the main sheet ---
# Each populated H cell advances to the next CU code, and the last one
# (row 20) gets a brand-new "CU11" code.
$ws = $wb.Worksheets.Item("Requerimientos consolidado")

$ws.Range("H8").Value  = "CU03"
$ws.Range("H10").Value = "CU04"
$ws.Range("H11").Value = "CU05"
$ws.Range("H13").Value = "CU06"
$ws.Range("H14").Value = "CU07"
$ws.Range("H16").Value = "CU08"
$ws.Range("H17").Value = "CU09"
$ws.Range("H18").Value = "CU10"
$ws.Range("H20").Value = "CU11"

# --- Update the sheet view: zoom level and active-cell selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
[void]$ws.Range("H21").Select()

# --- Remove the leftover "Hoja1" worksheet ---
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Hoja1").Delete()
$excel.DisplayAlerts = $true
